# Updates cryptos list values (Price / Volume(1h)) and one coin row (49 -> Mantle)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new Price (D) text value (or $null to leave as-is), new Volume(1h) (E) text value
$updates = @(
    @(2, "87.548.19", "  +6.58%  "),
    @(3, "3.294.40", "  +3.11%  "),
    @(4, $null, "  +0.32%  "),
    @(5, "213.43", "  -1.29%  "),
    @(6, "628.91", "  +0.64%  "),
    @(7, "0.396", "  +35.17%  "),
    @(8, "0.648", "  +10.56%  "),
    @(9, $null, "  +0.14%  "),
    @(10, "3.293.34", "  +3.18%  "),
    @(11, "0.591", "  -0.28%  "),
    @(12, "0.0000266", "  +2.01%  "),
    @(13, $null, "  +7.21%  "),
    @(14, "34.62", "  +8.01%  "),
    @(15, "3.909.03", "  +3.42%  "),
    @(16, "5.30", "  -0.54%  "),
    @(17, "87.418.46", "  +6.70%  "),
    @(18, "3.296.29", "  +3.40%  "),
    @(19, "14.31", "  +1.42%  "),
    @(20, "3.01", "  -7.39%  "),
    @(21, "9.24", "  +2.40%  "),
    @(22, "437.93", "  +0.33%  "),
    @(23, "5.40", "  +4.72%  "),
    @(24, "7.25", "  -0.29%  "),
    @(25, "12.37", "  +9.75%  "),
    @(26, "5.24", "  -3.04%  "),
    @(27, $null, "  +2.96%  "),
    @(28, "77.26", "  +0.69%  "),
    @(29, "0.0000132", "  +6.19%  "),
    @(30, "0.998", "  +0.05%  "),
    @(31, "0.183", "  +24.77%  "),
    @(32, $null, "  +1.24%  "),
    @(33, "9.02", "  -1.14%  "),
    @(34, "555.79", "  -6.08%  "),
    @(35, "1.46", "  -4.34%  "),
    @(36, "1.98", "  -1.63%  "),
    @(37, "6.99", "  +13.17%  "),
    @(38, "0.140", "  -9.63%  "),
    @(39, "22.80", "  -0.42%  "),
    @(40, $null, "  +0.23%  "),
    @(41, "21.77", "  +4.62%  "),
    @(42, "0.403", "  -1.83%  "),
    @(43, "2.04", "  -1.56%  "),
    @(44, "2.99", "  -3.10%  "),
    @(45, $null, "  +0.02%  "),
    @(46, "155.35", "  -3.54%  "),
    @(47, "181.68", "  -3.57%  "),
    @(48, "1.36", "  +1.10%  "),
    @(49, "45.20", "  +1.07%  "),
    @(50, "4.29", "  +0.93%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceVal = $u[1]
    $volVal = $u[2]
    if ($null -ne $priceVal) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $priceVal
        $dCell.Style = "Normal"
    }
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $volVal
    $eCell.Style = "Normal"
}

# Row 51: coin changed from ARBITRUM to Mantle
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$d51 = $ws.Cells.Item(51, 4)
$d51.NumberFormat = "@"
$d51.Value = "0.753"
$d51.Style = "Normal"
$e51 = $ws.Cells.Item(51, 5)
$e51.NumberFormat = "@"
$e51.Value = "  -2.80%  "
$e51.Style = "Normal"
